$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 139, shifting existing rows 139:149 down to 140:150.
$ws.Rows("139:139").Insert()

# Populate the newly inserted row 139 with the new record.
$ws.Range("A139").Value = 7
$ws.Range("B139").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C139").Value = "Ñuble"
$ws.Range("D139").Value = 44516
$ws.Range("D139").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E139").Value = 16
$ws.Range("F139").Value = 100112006
$ws.Range("G139").Value = "Repollo"
$ws.Range("H139").Value = "Crespo record"
$ws.Range("I139").Value = "Primera"
$ws.Range("J139").Value = 300
$ws.Range("K139").Value = 600
$ws.Range("L139").Value = 700
$ws.Range("M139").Value = 650
$ws.Range("N139").Value = "$/unidad"
$ws.Range("O139").Value = "Provincia de Diguillín"
$ws.Range("P139").Value = 650
$ws.Range("Q139").Value = 1
$ws.Range("R139").Value = "Hortaliza"
